# "added test for symmetric power net"
#
# The timesheet gets one more logged interval:
#  - the existing (until now unfinished) row 178 entry is completed
#    (its day changes from the 8th to the 7th, and an end time /
#    duration formulas are filled in)
#  - a brand new row is inserted right after it for another entry on
#    the 7th (start time only, end time left blank like before)
#  - the trailing "blank separator" row and the three summary rows
#    (sum [min] / sum [h] / sum [working weeks]) shift down by one row
#    as a consequence of the insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 179 - this pushes the blank
# separator row and the three summary rows down by one (179->180,
# 180->181, 181->182, 182->183) and grows the used range to H183.
$ws.Rows(179).Insert()

# Row 178 was a half-entered entry (only a start time). Finish it:
# correct the day, and add the end time + the two duration formulas
# that mirror the ones used by every other row in this block.
$ws.Range("C178").Value = 7
$ws.Range("E178").Value = 0.625
$ws.Range("F178").Formula = "=(E178-D178)*24*60"
$ws.Range("G178").Formula = "=F178/60"

# New row 179: another entry on the same day, start time only (end
# time / duration are left blank, same as row 178 used to be).
$ws.Range("A179").Value = 2014
$ws.Range("B179").Value = 8
$ws.Range("C179").Value = 7
$ws.Range("D179").Value = 0.65277777777777779

# Move the selection onto the newly-blank end-time cell.
$ws.Range("E179").Select()
